$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-12-12 Friday"; new="2025-12-13 Saturday"},
    @{old="56×31="; new="48×86="},
    @{old="70×11="; new="55×48="},
    @{old="73×32="; new="68×53="},
    @{old="89×26="; new="45×30="},
    @{old="83×31="; new="58×19="},
    @{old="20×11="; new="90×64="},
    @{old="75×49="; new="63×89="},
    @{old="16×13="; new="48×57="},
    @{old="79×67="; new="39×87="},
    @{old="56×72="; new="17×44="},
    @{old="66×29="; new="51×12="},
    @{old="16×11="; new="35×68="},
    @{old="17×62="; new="54×52="},
    @{old="88×13="; new="87×65="},
    @{old="19×34="; new="75×97="},
    @{old="26×57="; new="78×29="},
    @{old="54×54="; new="43×67="},
    @{old="30×16="; new="66×42="},
    @{old="47×33="; new="65×90="},
    @{old="16×96="; new="25×94="},
    @{old="39×88="; new="88×48="},
    @{old="13×71="; new="93×66="},
    @{old="56×49="; new="76×36="},
    @{old="16×39="; new="18×41="},
    @{old="48×58="; new="85×22="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
